# Generate Report for Handoff
#
# The localization-status report is regenerated: the two tracked files
# (61b896cf-...md and e0c3cbec-...md) swap row order in every sheet, and
# 61b896cf-...md picks up a new handoff event ("Ready for handoff" /
# updated handoff timestamps) while e0c3cbec-...md keeps its previous
# "Handed back: in sync with en-US" status.

$wb = $excel.ActiveWorkbook

$file61 = "61b896cf-cc6b-4613-bae6-25589e9c641c.md"
$file61zhxlf = "61b896cf-cc6b-4613-bae6-25589e9c641c.0cb423db10d2ca3cac4e4e2e5696829bdf7b154d.zh-cn.xlf"
$file61dexlf = "61b896cf-cc6b-4613-bae6-25589e9c641c.0cb423db10d2ca3cac4e4e2e5696829bdf7b154d.de-de.xlf"

$fileE0 = "e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.md"
$fileE0zhxlf = "e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.f46fd9bbdb5bce68e26b2f9491a78b463d29c64c.zh-cn.xlf"
$fileE0dexlf = "e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.f46fd9bbdb5bce68e26b2f9491a78b463d29c64c.de-de.xlf"

$statusHandedBack = "Handed back: in sync with en-US"
$statusReady = "Ready for handoff"

# Original (still valid) external link targets -- unchanged from before.xlsx,
# only which row/display uses them changes.
$url61md    = "https://github.com/OpenLocalizationTest/oltest/blob/822ada9a4bcf95cb0233b57cd112bbaf790e1904/e2e/61b896cf-cc6b-4613-bae6-25589e9c641c.md"
$urlE0md    = "https://github.com/OpenLocalizationTest/oltest/blob/822ada9a4bcf95cb0233b57cd112bbaf790e1904/e2e/e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.md"

$url61zhxlf_off = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/aadc418136b2ab09f0d6912698843f0a83ab67e5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/61b896cf-cc6b-4613-bae6-25589e9c641c.0cb423db10d2ca3cac4e4e2e5696829bdf7b154d.zh-cn.xlf"
$url61md_zh      = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/0e70a2396f99756c5b961090ec78f23f5e031637/e2e/61b896cf-cc6b-4613-bae6-25589e9c641c.md"
$url61zhxlf_back = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/db25dd1e22eccbf40b3b1a0c76d66d11b51ece47/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/61b896cf-cc6b-4613-bae6-25589e9c641c.0cb423db10d2ca3cac4e4e2e5696829bdf7b154d.zh-cn.xlf"

$urlE0zhxlf_off = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/aadc418136b2ab09f0d6912698843f0a83ab67e5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.f46fd9bbdb5bce68e26b2f9491a78b463d29c64c.zh-cn.xlf"
$urlE0md_zh      = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/0e70a2396f99756c5b961090ec78f23f5e031637/e2e/e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.md"
$urlE0zhxlf_back = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/db25dd1e22eccbf40b3b1a0c76d66d11b51ece47/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.f46fd9bbdb5bce68e26b2f9491a78b463d29c64c.zh-cn.xlf"

$url61dexlf_off = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b2e6fc8ab0cd39017e653a934ebb2e062a927978/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/61b896cf-cc6b-4613-bae6-25589e9c641c.0cb423db10d2ca3cac4e4e2e5696829bdf7b154d.de-de.xlf"
$url61md_de      = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e1e1db449caf8965ad2c0aeebd4be28195f61839/e2e/61b896cf-cc6b-4613-bae6-25589e9c641c.md"
$url61dexlf_back = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/664289f9bae1facc41161c1b4ab97b5c5337cc18/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/61b896cf-cc6b-4613-bae6-25589e9c641c.0cb423db10d2ca3cac4e4e2e5696829bdf7b154d.de-de.xlf"

$urlE0dexlf_off = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b2e6fc8ab0cd39017e653a934ebb2e062a927978/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.f46fd9bbdb5bce68e26b2f9491a78b463d29c64c.de-de.xlf"
$urlE0md_de      = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e1e1db449caf8965ad2c0aeebd4be28195f61839/e2e/e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.md"
$urlE0dexlf_back = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/664289f9bae1facc41161c1b4ab97b5c5337cc18/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.f46fd9bbdb5bce68e26b2f9491a78b463d29c64c.de-de.xlf"

# ----------------------------------------------------------------------
# Overview sheet: row 2 becomes e0c3cbec (unchanged status/date), row 3
# becomes 61b896cf with the new "Ready for handoff" status and the new
# handoff date.
# ----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $fileE0
$wsOverview.Range("B2").Value = $statusHandedBack
$wsOverview.Range("C2").Value = $statusHandedBack
$wsOverview.Range("D2").Value = "2016-03-23 22:48:57"

$wsOverview.Range("A3").Value = $file61
$wsOverview.Range("B3").Value = $statusReady
$wsOverview.Range("C3").Value = $statusReady
$wsOverview.Range("D3").Value = "2016-03-23 22:50:31"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $url61md, "", "", $fileE0)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $urlE0md, "", "", $file61)

# ----------------------------------------------------------------------
# zh-cn sheet: same row swap; handoff file/datetime for 61b896cf updated
# to the new xlf handoff.
# ----------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $fileE0
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = $statusHandedBack
$wsZh.Range("D2").Value = $fileE0zhxlf
$wsZh.Range("E2").Value = "2016-03-23 22:48:53"
$wsZh.Range("F2").Value = $fileE0
$wsZh.Range("G2").Value = $fileE0zhxlf
$wsZh.Range("H2").Value = "2016-03-23 22:49:31"
$wsZh.Range("J2").Value = "Include"

$wsZh.Range("A3").Value = $file61
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = $statusReady
$wsZh.Range("D3").Value = $file61zhxlf
$wsZh.Range("E3").Value = "2016-03-23 22:50:27"
$wsZh.Range("F3").Value = $file61
$wsZh.Range("G3").Value = $file61zhxlf
$wsZh.Range("H3").Value = "2016-03-23 22:49:31"
$wsZh.Range("J3").Value = "Include"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $url61md, "", "", $fileE0)
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $url61zhxlf_off, "", "", $fileE0zhxlf)
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $url61md_zh, "", "", $fileE0)
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $url61zhxlf_back, "", "", $fileE0zhxlf)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $urlE0md, "", "", $file61)
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $urlE0zhxlf_off, "", "", $file61zhxlf)
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $urlE0md_zh, "", "", $file61)
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $urlE0zhxlf_back, "", "", $file61zhxlf)

# ----------------------------------------------------------------------
# de-de sheet: same pattern as zh-cn.
# ----------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $fileE0
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = $statusHandedBack
$wsDe.Range("D2").Value = $fileE0dexlf
$wsDe.Range("E2").Value = "2016-03-23 22:48:57"
$wsDe.Range("F2").Value = $fileE0
$wsDe.Range("G2").Value = $fileE0dexlf
$wsDe.Range("H2").Value = "2016-03-23 22:49:40"
$wsDe.Range("J2").Value = "Include"

$wsDe.Range("A3").Value = $file61
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = $statusReady
$wsDe.Range("D3").Value = $file61dexlf
$wsDe.Range("E3").Value = "2016-03-23 22:50:31"
$wsDe.Range("F3").Value = $file61
$wsDe.Range("G3").Value = $file61dexlf
$wsDe.Range("H3").Value = "2016-03-23 22:49:40"
$wsDe.Range("J3").Value = "Include"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $url61md, "", "", $fileE0)
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $url61dexlf_off, "", "", $fileE0dexlf)
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $url61md_de, "", "", $fileE0)
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $url61dexlf_back, "", "", $fileE0dexlf)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $urlE0md, "", "", $file61)
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $urlE0dexlf_off, "", "", $file61dexlf)
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $urlE0md_de, "", "", $file61)
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $urlE0dexlf_back, "", "", $file61dexlf)
